# Rename the inline logo pictures living in the document's headers/footers.
#
# Mapping of Word's Headers/Footers COM collections to the underlying OOXML
# parts (confirmed via w:headerReference/w:footerReference w:type attrs):
#   Headers.Item(1) / Footers.Item(1) -> wdHeaderFooterPrimary ("default")
#   Headers.Item(2) / Footers.Item(2) -> wdHeaderFooterFirstPage ("first")
#
# Pearson logo (footers): image1.png -> image2.png
# BTEC logo   (headers): image2.jpg -> image1.jpg

$d = $word.ActiveDocument
$section = $d.Sections.Item(1)

# Footers - Pearson logo, rename image1.png -> image2.png
$footerPrimary = $section.Footers.Item(1)
$footerPrimary.Range.InlineShapes.Item(1).Name = "image2.png"

$footerFirst = $section.Footers.Item(2)
$footerFirst.Range.InlineShapes.Item(1).Name = "image2.png"

# Headers - BTEC logo, rename image2.jpg -> image1.jpg
$headerPrimary = $section.Headers.Item(1)
$headerPrimary.Range.InlineShapes.Item(1).Name = "image1.jpg"

$headerFirst = $section.Headers.Item(2)
$headerFirst.Range.InlineShapes.Item(1).Name = "image1.jpg"
